$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110, pushing the existing rows 110-135 down to 111-136.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly price record.
$ws.Cells.Item(110, 1).Value = 4
$ws.Cells.Item(110, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(110, 3).Value = "Los Lagos"
$ws.Cells.Item(110, 4).Value = 45027
$ws.Cells.Item(110, 5).Value = 10
$ws.Cells.Item(110, 6).Value = "Fruta"
$ws.Cells.Item(110, 7).Value = 100104
$ws.Cells.Item(110, 8).Value = "Frutos de pepita"
$ws.Cells.Item(110, 9).Value = 100104003
$ws.Cells.Item(110, 10).Value = "Membrillo"
$ws.Cells.Item(110, 11).Value = "Champion"
$ws.Cells.Item(110, 12).Value = "Primera"
$ws.Cells.Item(110, 13).Value = 300
$ws.Cells.Item(110, 14).Value = 15000
$ws.Cells.Item(110, 15).Value = 16000
$ws.Cells.Item(110, 16).Value = 15500
$ws.Cells.Item(110, 17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(110, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(110, 19).Value = 861
$ws.Cells.Item(110, 20).Value = 18
